$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.6126695871353149
$ws.Range("B1").Value = 1.337053537368774
$ws.Range("C1").Value = 2.470052719116211
$ws.Range("D1").Value = 5.866518020629883
$ws.Range("E1").Value = 2.077924489974976
